$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("J2").Value = 7085
$ws.Range("J3").Value = 7480
$ws.Range("J4").Value = 1632
$ws.Range("J6").Value = 10138
$ws.Range("J7").Value = 26919

$ws = $wb.Worksheets.Item("Grant Park")
$ws.Range("J5").Value = 5
$ws.Range("J6").Value = 14

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("J2").Value = 446
$ws.Range("J3").Value = 498
$ws.Range("J6").Value = 608
$ws.Range("J7").Value = 1683

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("J2").Value = 157
$ws.Range("J3").Value = 201
$ws.Range("J7").Value = 540

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("J3").Value = 404
$ws.Range("J6").Value = 435
$ws.Range("J7").Value = 1222

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("J3").Value = 139
$ws.Range("J7").Value = 388

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("J3").Value = 280
$ws.Range("J7").Value = 828

$ws = $wb.Worksheets.Item("New City")
$ws.Range("J6").Value = 251
$ws.Range("J7").Value = 677

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("J6").Value = 108
$ws.Range("J7").Value = 414

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("J2").Value = 215
$ws.Range("J4").Value = 123
$ws.Range("J6").Value = 206
$ws.Range("J7").Value = 769
$ws.Range("J8").Value = 1683
$ws.Range("J10").Value = 198
$ws.Range("J11").Value = 477
$ws.Range("J14").Value = 145
$ws.Range("J15").Value = 334
$ws.Range("J16").Value = 105
$ws.Range("J18").Value = 221
$ws.Range("J19").Value = 778
$ws.Range("J20").Value = 567
$ws.Range("J24").Value = 87
$ws.Range("J25").Value = 137
$ws.Range("J27").Value = 164
$ws.Range("J29").Value = 1440
$ws.Range("J31").Value = 275
$ws.Range("J33").Value = 1222
$ws.Range("J37").Value = 828
$ws.Range("J38").Value = 14
$ws.Range("J41").Value = 201
$ws.Range("J42").Value = 1158
$ws.Range("J47").Value = 197
$ws.Range("J50").Value = 161
$ws.Range("J51").Value = 332
$ws.Range("J52").Value = 686
$ws.Range("J54").Value = 526
$ws.Range("J55").Value = 427
$ws.Range("J60").Value = 157
$ws.Range("J63").Value = 86
$ws.Range("J64").Value = 178
$ws.Range("J65").Value = 677
$ws.Range("J67").Value = 1002
$ws.Range("J73").Value = 259
$ws.Range("J76").Value = 382
$ws.Range("J77").Value = 187
$ws.Range("J79").Value = 747
$ws.Range("J80").Value = 47
$ws.Range("J83").Value = 540
$ws.Range("J85").Value = 1108
$ws.Range("J87").Value = 87
$ws.Range("J89").Value = 336
$ws.Range("J93").Value = 114
$ws.Range("J95").Value = 388
$ws.Range("J96").Value = 294
$ws.Range("J99").Value = 414
$ws.Range("J101").Value = 26919

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("J2").Value = 94
$ws.Range("J7").Value = 275

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("J3").Value = 375
$ws.Range("J6").Value = 278
$ws.Range("J7").Value = 1002

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("J2").Value = 131
$ws.Range("J5").Value = 5
$ws.Range("J6").Value = 244
$ws.Range("J7").Value = 526

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("J2").Value = 435
$ws.Range("J3").Value = 509
$ws.Range("J6").Value = 366
$ws.Range("J7").Value = 1440

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("J2").Value = 189
$ws.Range("J3").Value = 222
$ws.Range("J4").Value = 37
$ws.Range("J6").Value = 300
$ws.Range("J7").Value = 778

$ws = $wb.Worksheets.Item("River North")
$ws.Range("J2").Value = 68
$ws.Range("J7").Value = 382

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("J4").Value = 9
$ws.Range("J7").Value = 145

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("J2").Value = 62
$ws.Range("J7").Value = 206

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("J6").Value = 124
$ws.Range("J7").Value = 201

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("J6").Value = 618
$ws.Range("J7").Value = 1158

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("J6").Value = 112
$ws.Range("J7").Value = 198

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("J6").Value = 242
$ws.Range("J7").Value = 427

$ws = $wb.Worksheets.Item("Dunning")
$ws.Range("J6").Value = 22
$ws.Range("J7").Value = 87

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("J3").Value = 78
$ws.Range("J6").Value = 106
$ws.Range("J7").Value = 294

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("J2").Value = 210
$ws.Range("J3").Value = 249
$ws.Range("J6").Value = 224
$ws.Range("J7").Value = 747

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("J3").Value = 44
$ws.Range("J4").Value = 19
$ws.Range("J7").Value = 178

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("J2").Value = 159
$ws.Range("J6").Value = 163
$ws.Range("J7").Value = 567

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("J3").Value = 46
$ws.Range("J7").Value = 221

$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("J3").Value = 36
$ws.Range("J7").Value = 114

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("J2").Value = 240
$ws.Range("J3").Value = 233
$ws.Range("J7").Value = 769

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("J2").Value = 55
$ws.Range("J6").Value = 28
$ws.Range("J7").Value = 137

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("J6").Value = 92
$ws.Range("J7").Value = 197

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("J2").Value = 94
$ws.Range("J6").Value = 152
$ws.Range("J7").Value = 334

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("J6").Value = 56
$ws.Range("J7").Value = 161

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("J2").Value = 135
$ws.Range("J3").Value = 83
$ws.Range("J6").Value = 224
$ws.Range("J7").Value = 477

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("J5").Value = 3
$ws.Range("J7").Value = 259

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("J2").Value = 65
$ws.Range("J7").Value = 215

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("J6").Value = 103
$ws.Range("J7").Value = 336

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("J3").Value = 40
$ws.Range("J7").Value = 164

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("J2").Value = 73
$ws.Range("J7").Value = 332

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("J3").Value = 44
$ws.Range("J7").Value = 157

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("J3").Value = 401
$ws.Range("J7").Value = 1108

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("J2").Value = 69
$ws.Range("J7").Value = 187

$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Range("J6").Value = 25
$ws.Range("J7").Value = 47

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("J2").Value = 161
$ws.Range("J6").Value = 293
$ws.Range("J7").Value = 686

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("J6").Value = 56
$ws.Range("J7").Value = 123

$ws = $wb.Worksheets.Item("Ukrainian Village")
$ws.Range("J2").Value = 11
$ws.Range("J7").Value = 87

$ws = $wb.Worksheets.Item("Bucktown")
$ws.Range("J6").Value = 81
$ws.Range("J7").Value = 105
